# Generate Report for Handoff
# The 495132bb... file is now "Ready for handoff" with a fresh handoff
# timestamp (2016-*-17 20:38:*). The b7e89ace... file's row is removed
# from every sheet (it no longer appears in this handoff batch).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-38-17 20:38:32"
$ws.Rows.Item(3).Delete()

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-17 20:38:23"
$ws.Rows.Item(3).Delete()

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-17 20:38:32"
$ws.Rows.Item(3).Delete()
